# Rename workbook sheets and update the corresponding label cells (A1)
# to match the new naming scheme.
#
#   pw_summary   -> purchase_summary_sheet   (cell A1: purchases_summary_list -> purchase_summary_list)
#   pw_purchases -> purchase_sheet           (cell A1: purchases_table        -> purchase_table)

$wb = $excel.ActiveWorkbook

$summarySheet  = $wb.Worksheets.Item("pw_summary")
$purchaseSheet = $wb.Worksheets.Item("pw_purchases")

$summarySheet.Range("A1").Value  = "purchase_summary_list"
$purchaseSheet.Range("A1").Value = "purchase_table"

$summarySheet.Name  = "purchase_summary_sheet"
$purchaseSheet.Name = "purchase_sheet"
